# encode_pairs.xlsx - "Fixed Trigger & Jitter"
# - selected/added jitter ITI values for the second block of scenes (rows 62-91)
# - removed the now-redundant column O (duplicate of column N)
# - removed trailing blank rows 92:94
# - minor view-state updates (zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the jitter_ITI values for rows 62:79 (previously empty) ---
$jitterBlock1 = @(2, 2.05, 2.1, 2.15, 2.2, 2.25, 2.3, 2.35, 2.4, 2.45, 2.5, 2.55, 2.6, 2.65, 2.7, 2.75, 2.8, 2.85)
for ($i = 0; $i -lt $jitterBlock1.Length; $i++) {
    $ws.Cells.Item(62 + $i, 14).Value = $jitterBlock1[$i]
}

# --- 2. Fill in the jitter_ITI values for rows 80:91 (previously blank placeholder
#        cells formatted with the "0.00" number style) - reset to the default
#        General style to match the rest of the column, then set the value ---
$jitterBlock2 = @(2.9, 2.95, 3, 3.05, 3.1, 3.15, 3.2, 3.25, 2, 2.05, 2.1, 2.85)
for ($i = 0; $i -lt $jitterBlock2.Length; $i++) {
    $cell = $ws.Cells.Item(80 + $i, 14)
    $cell.Style = "Normal"
    $cell.Value = $jitterBlock2[$i]
}

# --- 3. Remove the trailing blank rows 92:94 ---
$ws.Range("A92:A94").EntireRow.Delete()

# --- 4. Remove column O (duplicate of the jitter_ITI column N) ---
$ws.Columns("O").Delete()

# --- 5. Update the view state: zoom in and move the selection down to row 61 ---
[void]$ws.Activate()
$excel.ActiveWindow.Zoom = 110
[void]$ws.Rows("61:61").EntireRow.Select()
